# Update the "想去人数" (interested-people count) figures in column F
# on the 展览, 演出, and 全部类型 worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 64
$ws1.Range("F5").Value  = 506
$ws1.Range("F6").Value  = 1498
$ws1.Range("F7").Value  = 984
$ws1.Range("F8").Value  = 111
$ws1.Range("F9").Value  = 209
$ws1.Range("F10").Value = 148
$ws1.Range("F11").Value = 208
$ws1.Range("F13").Value = 180
$ws1.Range("F14").Value = 165

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 3

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 64
$ws4.Range("F5").Value  = 506
$ws4.Range("F6").Value  = 1498
$ws4.Range("F7").Value  = 3
$ws4.Range("F8").Value  = 984
$ws4.Range("F9").Value  = 111
$ws4.Range("F10").Value = 209
$ws4.Range("F11").Value = 148
$ws4.Range("F12").Value = 208
$ws4.Range("F14").Value = 180
$ws4.Range("F15").Value = 165
